$d = $word.ActiveDocument

# Paragraph 22 is the empty paragraph right after the "POODLE Index..." paragraph,
# and it is directly followed (originally) by "Printing from crawled...". We turn it
# into a new bold+underlined "TODO" note and move the _GoBack bookmark from the old
# "Depth 2..." paragraph onto the end of this new paragraph.

# 1) Grab the (currently empty) paragraph and anchor a bookmark on it while it is
#    still empty -- this keeps the bookmark collapsed at a position that will end up
#    right after the paragraph's only run once we insert text before it.
$d.Bookmarks.Item("_GoBack").Delete()
$p22 = $d.Paragraphs.Item(22)
$p22.Range.Bookmarks.Add("_GoBack")

# 2) Insert the note text before the (collapsed) bookmark, so the new run ends up
#    before the bookmark markers, matching the original bookmark's placement pattern.
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range
$bmRange.InsertBefore("TODO: Recusion code for search – highlight it")

# 3) Bold + underline the whole paragraph (text run and paragraph mark).
$newPara = $d.Paragraphs.Item(22)
$newPara.Range.Bold = 1
$newPara.Range.Underline = 1
